$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that sat between "MP73010"
#    and " - Assignment 1 exercise" in the title paragraph.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Collapse the ">>  " / ">  your" / " stuff after this line >>>"
#    three-run paragraph (with the gramStart/gramEnd proofing marks)
#    down into a single run with the same combined text.
# ------------------------------------------------------------------
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq ">>>  your stuff after this line >>>") {
        $targetPara = $p
        break
    }
}
if ($targetPara -ne $null) {
    $targetPara.Range.Find.Execute(">>>  your stuff after this line >>>", $true, $false, $false, $false, $false, $true, 1, $false, ">>>  your stuff after this line >>>", 2)
}

# ------------------------------------------------------------------
# 3) Replace the first trailing empty paragraph (directly after
#    "Ben changing things up!") with the new commentary paragraphs,
#    leaving the final empty paragraph at the very end untouched.
# ------------------------------------------------------------------
$benPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Ben changing things up!") {
        $benPara = $p
        break
    }
}

$insertAfter = $benPara.Range
$insertAfter.Collapse(0)
$insertAfter.InsertParagraphAfter()

$p1 = $benPara.Next()
$p1.Range.Text = "Is ben changing it up or has Version control been around for a long time, what exactly is ben changing up? without context it" + [char]8217 + "s a rather confusing statement."

$p1.Range.InsertParagraphAfter()
$p2 = $p1.Next()
$p2.Range.Text = "Version control systems are a critical element to all software development,"

$lb = [char]11
$p2r = $p2.Range
$insPos = $d.Range($p2r.End - 1, $p2r.End - 1)
$insPos.InsertAfter($lb + "The benefits are ranged, the immediate benefits that comes to mind is as the title suggests Version Control, the capacity to traverse previous versions in the event new code breaks something its easy to compare and see what is different.")

$p2r = $p2.Range
$insPos = $d.Range($p2r.End - 1, $p2r.End - 1)
$insPos.InsertAfter($lb)

$p2r = $p2.Range
$insPos = $d.Range($p2r.End - 1, $p2r.End - 1)
$insPos.InsertAfter($lb + "Additionally version control ensures can be used as a backup solution provided the VCS platform is itself backed up, it" + [char]8217 + "s a rather robust solution, most platforms are offered as cloud services with inherent backup and redundancy features.")

$p2.Range.InsertParagraphAfter()
$p3 = $p2.Next()
$p3.Range.Text = "Version control systems are also absolutely required for open source to progress."

# ------------------------------------------------------------------
# 4) Re-add the "_GoBack" bookmark, now spanning (collapsed) at the
#    end of the new final paragraph's text.
# ------------------------------------------------------------------
$p3r = $p3.Range
$bmRange = $d.Range($p3r.End - 1, $p3r.End - 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
